# Lattice multiplication exercises worksheet: regenerate the practice
# problems (the header "A x B" plus the partial-product / lattice grid
# scaffolding) for every cell of the 5x3 table.
#
# Word's Find/Replace on a Range in this host always matches against the
# whole document instead of being scoped to the Range it was invoked on,
# and several of the per-cell text fragments (e.g. "  4    3") repeat
# across different cells -- so naive Find/Replace would clobber the wrong
# cell. Instead we set each table cell's Range.Text directly, using
# Chr(11) (the vertical-tab code Word uses internally for a manual line
# break / <w:br/>) to join the five lines of each cell so the line-break
# structure of the cell is preserved.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

# Row, Col, [header, partials, rule, left-carry, right-carry]
$cells = @(
    @(1, 1, @("19 x 57", "  5    7", "  ----", "1|    |", "9|    |")),
    @(1, 2, @("93 x 60", "  6    0", "  ----", "9|    |", "3|    |")),
    @(1, 3, @("19 x 37", "  3    7", "  ----", "1|    |", "9|    |")),

    @(2, 1, @("41 x 88", "  8    8", "  ----", "4|    |", "1|    |")),
    @(2, 2, @("76 x 38", "  3    8", "  ----", "7|    |", "6|    |")),
    @(2, 3, @("65 x 44", "  4    4", "  ----", "6|    |", "5|    |")),

    @(3, 1, @("60 x 96", "  9    6", "  ----", "6|    |", "0|    |")),
    @(3, 2, @("49 x 63", "  6    3", "  ----", "4|    |", "9|    |")),
    @(3, 3, @("37 x 97", "  9    7", "  ----", "3|    |", "7|    |")),

    @(4, 1, @("51 x 22", "  2    2", "  ----", "5|    |", "1|    |")),
    @(4, 2, @("43 x 89", "  8    9", "  ----", "4|    |", "3|    |")),
    @(4, 3, @("91 x 22", "  2    2", "  ----", "9|    |", "1|    |")),

    @(5, 1, @("15 x 81", "  8    1", "  ----", "1|    |", "5|    |")),
    @(5, 2, @("68 x 72", "  7    2", "  ----", "6|    |", "8|    |")),
    @(5, 3, @("32 x 43", "  4    3", "  ----", "3|    |", "2|    |"))
)

foreach ($entry in $cells) {
    $row = $entry[0]
    $col = $entry[1]
    $lines = $entry[2]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = ($lines -join $nl)
}
